$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
try {
  $d2 = $nm.Design
  Write-Host "NotesMaster.Design = $d2 Name=$($d2.Name)"
} catch {
  Write-Host "err: $_"
}
